$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Sheet1" to "Bugs"
$ws1 = $wb.ActiveSheet
$ws1.Name = "Bugs"

# Add a new worksheet right after "Bugs" and name it "Improvements"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Improvements"

# Populate the "Improvements" sheet.
# Write the longer description text first so it lands at shared-string
# index 10, then the "Improvements" title so it lands at index 11.
$ws2.Range("C4").Value = "How not to replicate code between companion-creation.component and companion-management.component"
$ws2.Range("C2").Value = "Improvements"
$ws2.Range("B3").Value = "Number"
$ws2.Range("C3").Value = "Description"
$ws2.Range("B4").Value = 1

# Match page setup / margins used by the rest of the workbook
$ws2.PageSetup.LeftMargin = 54
$ws2.PageSetup.RightMargin = 54
$ws2.PageSetup.TopMargin = 72
$ws2.PageSetup.BottomMargin = 72
$ws2.PageSetup.HeaderMargin = 36
$ws2.PageSetup.FooterMargin = 36
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selection / active sheet: "Improvements" is the active (selected) tab,
# with B5 as the selected cell
$ws2.Range("B5").Select() | Out-Null
